$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 10:31"

# 2. Update numeric data for rows whose country ranking did not change

# Row 6 - India
$ws.Range("B6").Value = 3049855
$ws.Range("C6").Value = 6419
$ws.Range("D6").Value = 2281982
$ws.Range("E6").Value = 710998
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 29
$ws.Range("H6").Value = 56875

# Row 7 - Rusia
$ws.Range("B7").Value = 956749
$ws.Range("C7").Value = 4852
$ws.Range("D7").Value = 770639
$ws.Range("E7").Value = 169727
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 73
$ws.Range("H7").Value = 16383

# Row 26 - Indonesia
$ws.Range("B26").Value = 153535
$ws.Range("C26").Value = 2037
$ws.Range("D26").Value = 107500
$ws.Range("E26").Value = 39355
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 86
$ws.Range("H26").Value = 6680

# Row 54 - Barein
$ws.Range("E54").Value = 3265
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 184

# Row 72 - Australia
$ws.Range("D72").Value = 19233
$ws.Range("E72").Value = 5077

# Row 81 - Bulgaria
$ws.Range("B81").Value = 15227
$ws.Range("D81").Value = 10322
$ws.Range("E81").Value = 4360
$ws.Range("H81").Value = 545

# Row 109 - Hungria
$ws.Range("B109").Value = 5155
$ws.Range("C109").Value = 22
$ws.Range("D109").Value = 3695
$ws.Range("E109").Value = 847
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = 613

# Row 121 - Eslovaquia
$ws.Range("B121").Value = 3356
$ws.Range("C121").Value = 40
$ws.Range("D121").Value = 2148
$ws.Range("E121").Value = 1175

# Row 150 - Georgia
$ws.Range("B150").Value = 1411
$ws.Range("C150").Value = 17
$ws.Range("E150").Value = 262

# Row 160 - Vietnam
$ws.Range("E160").Value = 424
$ws.Range("G160").Value = 1
$ws.Range("H160").Value = 27

# 3. Ucrania overtakes Kazajistan in the ranking: row 31 becomes Ucrania
#    (with fresh data) and row 32 becomes Kazajistan (keeping the totals
#    that used to belong to row 31).
$ws.Range("A31").Value = "Ucrania"
$ws.Range("B31").Value = 104958
$ws.Range("C31").Value = 1987
$ws.Range("D31").Value = 52235
$ws.Range("E31").Value = 50452
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 27
$ws.Range("H31").Value = 2271

$ws.Range("A32").Value = "Kazajistan"
$ws.Range("B32").Value = 104543
$ws.Range("C32").Value = 230
$ws.Range("D32").Value = 91089
$ws.Range("E32").Value = 12039
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 1415
